# daily auto push: 2025-10-03 13:33 UTC
# Append the new day's ranking row (row 56) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 56

# Column A holds dates stored as plain text (e.g. "2025/09/22"), not real
# Excel dates. Force the cell to Text format before assigning so the
# "2025/10/03" string isn't auto-converted into a date serial value, then
# clear the formatting again so the cell keeps the same (default) style as
# all the other data rows.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025/10/03"
$ws.Range("A$newRow").ClearFormats()

$ws.Range("B$newRow").Value = "金"
$ws.Range("C$newRow").Value = 20
$ws.Range("D$newRow").Value = 4
